# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (column E, rows 16-21) was re-sorted from
# descending (2203 .. 2110) to ascending (2110 .. 2203) order, and the
# one "Valor Mora" figure (column F) that differed from the rest moved
# along with its period (it was on the 2203 row, now at the bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E16").Value = "2110"
$ws.Range("E17").Value = "2111"
$ws.Range("E18").Value = "2112"
$ws.Range("E19").Value = "2201"
$ws.Range("E20").Value = "2202"
$ws.Range("E21").Value = "2203"

$ws.Range("F16").Value = 36341
$ws.Range("F17").Value = 36341
$ws.Range("F18").Value = 36341
$ws.Range("F19").Value = 36341
$ws.Range("F20").Value = 36341
$ws.Range("F21").Value = 32707
